$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the crypto-price refresh diff.
# "ForceText" marks cells whose new value parses as a plain number
# (e.g. "1.00", "0.0691") so we must pin them to text storage --
# otherwise Excel auto-converts and we lose the original "69.53" style formatting.
$updates = @(
    @{ Cell = 'D2'; Value = '69.540.83'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  +3.17%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '3.370.80'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  +4.67%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '192.45'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  +5.34%  '; ForceText = $false },
    @{ Cell = 'D6'; Value = '592.92'; ForceText = $true },
    @{ Cell = 'E6'; Value = '  +2.74%  '; ForceText = $false },
    @{ Cell = 'B7'; Value = 'USDC'; ForceText = $false },
    @{ Cell = 'C7'; Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; ForceText = $false },
    @{ Cell = 'D7'; Value = '1.00'; ForceText = $true },
    @{ Cell = 'E7'; Value = '  +0.16%  '; ForceText = $false },
    @{ Cell = 'B8'; Value = 'XRP'; ForceText = $false },
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; ForceText = $false },
    @{ Cell = 'D8'; Value = '0.609'; ForceText = $true },
    @{ Cell = 'E8'; Value = '  +1.09%  '; ForceText = $false },
    @{ Cell = 'E9'; Value = '  +3.38%  '; ForceText = $false },
    @{ Cell = 'E10'; Value = '  +2.94%  '; ForceText = $false },
    @{ Cell = 'D11'; Value = '0.421'; ForceText = $true },
    @{ Cell = 'E11'; Value = '  +2.51%  '; ForceText = $false },
    @{ Cell = 'D12'; Value = '3.958.99'; ForceText = $false },
    @{ Cell = 'E12'; Value = '  +4.81%  '; ForceText = $false },
    @{ Cell = 'E13'; Value = '  +1.21%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '28.66'; ForceText = $true },
    @{ Cell = 'E14'; Value = '  +3.39%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '69.564.69'; ForceText = $false },
    @{ Cell = 'E15'; Value = '  +3.10%  '; ForceText = $false },
    @{ Cell = 'E16'; Value = '  +2.27%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '3.419.75'; ForceText = $false },
    @{ Cell = 'E17'; Value = '  +6.12%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '449.40'; ForceText = $true },
    @{ Cell = 'E18'; Value = '  +13.83%  '; ForceText = $false },
    @{ Cell = 'E19'; Value = '  +1.67%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '13.85'; ForceText = $true },
    @{ Cell = 'E20'; Value = '  +3.41%  '; ForceText = $false },
    @{ Cell = 'E21'; Value = '  +3.83%  '; ForceText = $false },
    @{ Cell = 'D22'; Value = '73.60'; ForceText = $true },
    @{ Cell = 'E22'; Value = '  +3.65%  '; ForceText = $false },
    @{ Cell = 'E23'; Value = '  +0.01%  '; ForceText = $false },
    @{ Cell = 'D24'; Value = '3.520.21'; ForceText = $false },
    @{ Cell = 'E24'; Value = '  +4.70%  '; ForceText = $false },
    @{ Cell = 'D25'; Value = '0.519'; ForceText = $true },
    @{ Cell = 'E25'; Value = '  +1.22%  '; ForceText = $false },
    @{ Cell = 'E26'; Value = '  +4.24%  '; ForceText = $false },
    @{ Cell = 'E27'; Value = '  +4.61%  '; ForceText = $false },
    @{ Cell = 'D28'; Value = '9.58'; ForceText = $true },
    @{ Cell = 'E28'; Value = '  +0.41%  '; ForceText = $false },
    @{ Cell = 'E29'; Value = '  +0.08%  '; ForceText = $false },
    @{ Cell = 'E30'; Value = '  +2.51%  '; ForceText = $false },
    @{ Cell = 'D32'; Value = '5.61'; ForceText = $true },
    @{ Cell = 'E32'; Value = '  +1.09%  '; ForceText = $false },
    @{ Cell = 'D33'; Value = '1.29'; ForceText = $true },
    @{ Cell = 'E33'; Value = '  +4.14%  '; ForceText = $false },
    @{ Cell = 'D34'; Value = '7.04'; ForceText = $true },
    @{ Cell = 'E34'; Value = '  +1.65%  '; ForceText = $false },
    @{ Cell = 'E35'; Value = '  +0.01%  '; ForceText = $false },
    @{ Cell = 'E36'; Value = '  +4.05%  '; ForceText = $false },
    @{ Cell = 'D37'; Value = '165.01'; ForceText = $true },
    @{ Cell = 'E37'; Value = '  +2.67%  '; ForceText = $false },
    @{ Cell = 'D38'; Value = '1.95'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  +4.13%  '; ForceText = $false },
    @{ Cell = 'D39'; Value = '27.36'; ForceText = $true },
    @{ Cell = 'E39'; Value = '  +4.52%  '; ForceText = $false },
    @{ Cell = 'D40'; Value = '0.821'; ForceText = $true },
    @{ Cell = 'E40'; Value = '  +2.36%  '; ForceText = $false },
    @{ Cell = 'E41'; Value = '  +1.28%  '; ForceText = $false },
    @{ Cell = 'D42'; Value = '6.53'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  +0.71%  '; ForceText = $false },
    @{ Cell = 'D43'; Value = '2.744.74'; ForceText = $false },
    @{ Cell = 'D44'; Value = '2.55'; ForceText = $true },
    @{ Cell = 'E44'; Value = '  +3.67%  '; ForceText = $false },
    @{ Cell = 'D45'; Value = '25.65'; ForceText = $true },
    @{ Cell = 'E45'; Value = '  +4.91%  '; ForceText = $false },
    @{ Cell = 'D46'; Value = '0.0691'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  +1.28%  '; ForceText = $false },
    @{ Cell = 'D47'; Value = '344.83'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  +3.52%  '; ForceText = $false },
    @{ Cell = 'D48'; Value = '40.83'; ForceText = $true },
    @{ Cell = 'E48'; Value = '  +0.85%  '; ForceText = $false },
    @{ Cell = 'D49'; Value = '0.0286'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  +3.36%  '; ForceText = $false },
    @{ Cell = 'B50'; Value = 'Arweave'; ForceText = $false },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'; ForceText = $false },
    @{ Cell = 'D50'; Value = '32.97'; ForceText = $true },
    @{ Cell = 'E50'; Value = '  +7.85%  '; ForceText = $false },
    @{ Cell = 'B51'; Value = 'ONDO'; ForceText = $false },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'; ForceText = $false },
    @{ Cell = 'D51'; Value = '1.04'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  +7.79%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Pin number-format to Text before assigning, then restore the
        # default "Normal" style so no stray formatting is left behind.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
